$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 119) holds a date serial number ("Förändrad" / last-changed date)
# that was bumped by one day (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12).
$ws.Range("C2:C119").Value = 45181
